# Applies the "fixed bugs in inventory and bugs in NDC lookup" edit to the
# Manufacturer Packaging Slip (Invoice Template) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice Template")

# ---------------------------------------------------------------------------
# Helper: force a literal text value into a cell even when the text looks
# like a number or a date (e.g. "1", "02/23/2024"), so Excel doesn't silently
# convert it to a numeric / date serial value.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------------------
# Header block: invoice date / invoice number
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("B5") "02/23/2024"
$ws.Range("L5").Value = "QR02232024AP301"

# ---------------------------------------------------------------------------
# Column headings above the account blocks
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "Ship To:"
$ws.Range("E8").Value = "Bill To:"

# ---------------------------------------------------------------------------
# Account / Wholesaler / Credit-to address blocks (rows 10-15)
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "AUTREY PHARMACY 3"
$ws.Range("D10").Value = "Lake Carmel PhaRxmacy"
$ws.Range("E10").Value = "AUTREY PHARMACY 3"

$ws.Range("B11").Value = "800 E. ALTON GLOR BLVD, UNIT B"
$ws.Range("D11").Value = "1205 Central Blvd,"
$ws.Range("E11").Value = "800 E. ALTON GLOR BLVD, UNIT B"

$ws.Range("B12").Value = "BROWNSVILLE, TX, 78526"
$ws.Range("D12").Value = "Brownsville, TX 78520"
$ws.Range("E12").Value = "BROWNSVILLE, TX, 78526"

$ws.Range("B14").Value = "Phone: 956-525-7759, fax: "
$ws.Range("D14").Value = "Phone: (956) 548-0801"
$ws.Range("E14").Value = "Phone: 956-525-7759, fax: "

$ws.Range("B15").Value = "DEA: FA5030010, Exp: 06/30/2023"
$ws.Range("D15").Value = "DEA: 1233934230"
$ws.Range("E15").Value = "DEA: FA5030010, Exp: 06/30/2023"

# ---------------------------------------------------------------------------
# Line item row 20: NDC lookup bug fix - Lot/Exp columns replaced with "1",
# and Full Qty / Price stored as text "1" instead of numeric 1.
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("G20") "1"
Set-TextValue $ws.Range("H20") "1"
Set-TextValue $ws.Range("J20") "1"
Set-TextValue $ws.Range("L20") "1"

# ---------------------------------------------------------------------------
# Remove the test/sample product rows 21-23 that were found to be bugs, and
# roll their estimated value into the Est Value column.
# ---------------------------------------------------------------------------
$ws.Range("B21:J21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("M21").Value = 175

$ws.Range("B22:J22").ClearContents()
$ws.Range("L22").ClearContents()

$ws.Range("B23:J23").ClearContents()
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 175

$ws.Range("M24").Value = 225

# ---------------------------------------------------------------------------
# Corrected grand total
# ---------------------------------------------------------------------------
$ws.Range("L34").Value = 1

# ---------------------------------------------------------------------------
# Selection moves from H10:K10 to B8:C8
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B8:C8").Select()
